$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 17.05853031862614
$ws.Cells.Item(2, 3).Value = 9.359861101404508
$ws.Cells.Item(2, 4).Value = 6.051771177070816
$ws.Cells.Item(2, 5).Value = 9.299332100212366
$ws.Cells.Item(2, 7).Value = 3.739354338117113
$ws.Cells.Item(2, 11).Value = 16.50047868558372
$ws.Cells.Item(2, 13).Value = 17.20843302999209
$ws.Cells.Item(2, 14).Value = 25.70560317131548

$ws.Cells.Item(3, 2).Value = 16.79876475369259
$ws.Cells.Item(3, 3).Value = 9.145377902845432
$ws.Cells.Item(3, 4).Value = 5.945417838199975
$ws.Cells.Item(3, 5).Value = 9.104316470554025
$ws.Cells.Item(3, 7).Value = 3.743486420456252
$ws.Cells.Item(3, 11).Value = 16.31850307376125
$ws.Cells.Item(3, 13).Value = 17.0530054347049
$ws.Cells.Item(3, 14).Value = 25.59884364571214

$ws.Cells.Item(4, 2).Value = 16.64265921272538
$ws.Cells.Item(4, 3).Value = 9.014401978392014
$ws.Cells.Item(4, 4).Value = 5.881135945390999
$ws.Cells.Item(4, 5).Value = 8.98577502026577
$ws.Cells.Item(4, 7).Value = 3.746151279829294
$ws.Cells.Item(4, 11).Value = 16.21069923562996
$ws.Cells.Item(4, 13).Value = 16.96193496298504
$ws.Cells.Item(4, 14).Value = 25.53403865635141

$ws.Cells.Item(5, 2).Value = 16.57998003147256
$ws.Cells.Item(5, 3).Value = 8.961293080089947
$ws.Cells.Item(5, 4).Value = 5.855233556602746
$ws.Cells.Item(5, 5).Value = 8.937842373318322
$ws.Cells.Item(5, 7).Value = 3.747269488176869
$ws.Cells.Item(5, 11).Value = 16.16780325408534
$ws.Cells.Item(5, 13).Value = 16.92595440777768
$ws.Cells.Item(5, 14).Value = 25.50782450302697

$ws.Cells.Item(6, 2).Value = 16.56963104576401
$ws.Cells.Item(6, 3).Value = 8.952492878602497
$ws.Cells.Item(6, 4).Value = 5.850951245327931
$ws.Cells.Item(6, 5).Value = 8.929907904062945
$ws.Cells.Item(6, 7).Value = 3.747457117970608
$ws.Cells.Item(6, 11).Value = 16.16074423738657
$ws.Cells.Item(6, 13).Value = 16.92004915400371
$ws.Cells.Item(6, 14).Value = 25.50348359861693

$ws.Cells.Item(7, 2).Value = 16.64181000595994
$ws.Cells.Item(7, 3).Value = 9.013684548238849
$ws.Cells.Item(7, 4).Value = 5.880785382151601
$ws.Cells.Item(7, 5).Value = 8.98512697362232
$ws.Cells.Item(7, 7).Value = 3.746166229597968
$ws.Cells.Item(7, 11).Value = 16.2101164749757
$ws.Cells.Item(7, 13).Value = 16.96144509208842
$ws.Cells.Item(7, 14).Value = 25.5336843272982

$ws.Cells.Item(8, 2).Value = 16.96830909030308
$ws.Cells.Item(8, 3).Value = 9.285808475534436
$ws.Cells.Item(8, 4).Value = 6.01491101622795
$ws.Cells.Item(8, 5).Value = 9.231885186048437
$ws.Cells.Item(8, 7).Value = 3.740752647837418
$ws.Cells.Item(8, 11).Value = 16.43694658128909
$ws.Cells.Item(8, 13).Value = 17.15395746164693
$ws.Cells.Item(8, 14).Value = 25.66864190074444

$ws.Cells.Item(9, 2).Value = 17.63180222960969
$ws.Cells.Item(9, 3).Value = 9.821528698868205
$ws.Cells.Item(9, 4).Value = 6.28440427778762
$ws.Cells.Item(9, 5).Value = 9.722148466704349
$ws.Cells.Item(9, 7).Value = 3.731144090128612
$ws.Cells.Item(9, 11).Value = 16.91078385191446
$ws.Cells.Item(9, 13).Value = 17.56446923597418
$ws.Cells.Item(9, 14).Value = 25.93890206268666

$ws.Cells.Item(10, 2).Value = 18.12847086893259
$ws.Cells.Item(10, 3).Value = 10.2117639827472
$ws.Cells.Item(10, 4).Value = 6.484203775544437
$ws.Cells.Item(10, 5).Value = 10.08213854310803
$ws.Cells.Item(10, 7).Value = 3.724690258793085
$ws.Cells.Item(10, 11).Value = 17.27360052530434
$ws.Cells.Item(10, 13).Value = 17.8838817310943
$ws.Cells.Item(10, 14).Value = 26.14045884747224

$ws.Cells.Item(11, 2).Value = 18.35537491963554
$ws.Cells.Item(11, 3).Value = 10.38765585448223
$ws.Cells.Item(11, 4).Value = 6.575049302385935
$ws.Cells.Item(11, 5).Value = 10.24504526829256
$ws.Cells.Item(11, 7).Value = 3.721883895293933
$ws.Cells.Item(11, 11).Value = 17.44118402934232
$ws.Cells.Item(11, 13).Value = 18.03253587884557
$ws.Cells.Item(11, 14).Value = 26.23270751843535

$ws.Cells.Item(12, 2).Value = 18.441351913655
$ws.Cells.Item(12, 3).Value = 10.45395770304382
$ws.Cells.Item(12, 4).Value = 6.609408645151616
$ws.Cells.Item(12, 5).Value = 10.30654711017473
$ws.Cells.Item(12, 7).Value = 3.720839680258875
$ws.Cells.Item(12, 11).Value = 17.50495233596501
$ws.Cells.Item(12, 13).Value = 18.08926399436895
$ws.Cells.Item(12, 14).Value = 26.26771140308979

$ws.Cells.Item(13, 2).Value = 18.4228341944435
$ws.Cells.Item(13, 3).Value = 10.43969303186484
$ws.Cells.Item(13, 4).Value = 6.602011158664078
$ws.Cells.Item(13, 5).Value = 10.29331091088568
$ws.Cells.Item(13, 7).Value = 3.721063750202421
$ws.Cells.Item(13, 11).Value = 17.49120588847283
$ws.Cells.Item(13, 13).Value = 18.07702792544998
$ws.Cells.Item(13, 14).Value = 26.26016966637762

$ws.Cells.Item(14, 2).Value = 18.3624477113548
$ws.Cells.Item(14, 3).Value = 10.39311703862878
$ws.Cells.Item(14, 4).Value = 6.57787710387086
$ws.Cells.Item(14, 5).Value = 10.2501091740769
$ws.Cells.Item(14, 7).Value = 3.721797617195544
$ws.Cells.Item(14, 11).Value = 17.44642443592225
$ws.Cells.Item(14, 13).Value = 18.03719446190546
$ws.Cells.Item(14, 14).Value = 26.23558593300235

$ws.Cells.Item(15, 2).Value = 18.32546369406791
$ws.Cells.Item(15, 3).Value = 10.36454622111838
$ws.Cells.Item(15, 4).Value = 6.563087814447166
$ws.Cells.Item(15, 5).Value = 10.22362060572097
$ws.Cells.Item(15, 7).Value = 3.72224953650379
$ws.Cells.Item(15, 11).Value = 17.41903294472464
$ws.Cells.Item(15, 13).Value = 18.01285072351398
$ws.Cells.Item(15, 14).Value = 26.2205366986144

$ws.Cells.Item(16, 2).Value = 18.11365423336483
$ws.Cells.Item(16, 3).Value = 10.20023014270562
$ws.Cells.Item(16, 4).Value = 6.478262846566009
$ws.Cells.Item(16, 5).Value = 10.07146941190019
$ws.Cells.Item(16, 7).Value = 3.724876256301222
$ws.Cells.Item(16, 11).Value = 17.26269474145002
$ws.Cells.Item(16, 13).Value = 17.87423036419606
$ws.Cells.Item(16, 14).Value = 26.13444066232163

$ws.Cells.Item(17, 2).Value = 17.98390194245119
$ws.Cells.Item(17, 3).Value = 10.09895893258639
$ws.Cells.Item(17, 4).Value = 6.426188414275844
$ws.Cells.Item(17, 5).Value = 9.977863666143159
$ws.Cells.Item(17, 7).Value = 3.7265207432951
$ws.Cells.Item(17, 11).Value = 17.16739560739503
$ws.Cells.Item(17, 13).Value = 17.79001724058291
$ws.Cells.Item(17, 14).Value = 26.08176040331335

$ws.Cells.Item(18, 2).Value = 17.90936786872223
$ws.Cells.Item(18, 3).Value = 10.04056108867687
$ws.Cells.Item(18, 4).Value = 6.396234451942363
$ws.Cells.Item(18, 5).Value = 9.923947320048718
$ws.Cells.Item(18, 7).Value = 3.72747880723662
$ws.Cells.Item(18, 11).Value = 17.11282419864663
$ws.Cells.Item(18, 13).Value = 17.74189836839998
$ws.Cells.Item(18, 14).Value = 26.05151333498483

$ws.Cells.Item(19, 2).Value = 17.88415105226931
$ws.Cells.Item(19, 3).Value = 10.02076515968605
$ws.Cells.Item(19, 4).Value = 6.38609329559393
$ws.Cells.Item(19, 5).Value = 9.9056810065933
$ws.Cells.Item(19, 7).Value = 3.727805290308187
$ws.Cells.Item(19, 11).Value = 17.09439057181653
$ws.Cells.Item(19, 13).Value = 17.72566220178545
$ws.Cells.Item(19, 14).Value = 26.04128161297648

$ws.Cells.Item(20, 2).Value = 17.99770498693924
$ws.Cells.Item(20, 3).Value = 10.10975540546707
$ws.Cells.Item(20, 4).Value = 6.431732308467311
$ws.Cells.Item(20, 5).Value = 9.987836584646487
$ws.Cells.Item(20, 7).Value = 3.726344423292601
$ws.Cells.Item(20, 11).Value = 17.17751571399067
$ws.Cells.Item(20, 13).Value = 17.79894926073183
$ws.Cells.Item(20, 14).Value = 26.08736286385206

$ws.Cells.Item(21, 2).Value = 18.38018390777281
$ws.Cells.Item(21, 3).Value = 10.40680634136942
$ws.Cells.Item(21, 4).Value = 6.584967273338796
$ws.Cells.Item(21, 5).Value = 10.26280415332231
$ws.Cells.Item(21, 7).Value = 3.721581561789685
$ws.Cells.Item(21, 11).Value = 17.45956992363385
$ws.Cells.Item(21, 13).Value = 18.04888305146048
$ws.Cells.Item(21, 14).Value = 26.24280490990272

$ws.Cells.Item(22, 2).Value = 18.63042910251563
$ws.Cells.Item(22, 3).Value = 10.59914162074905
$ws.Cells.Item(22, 4).Value = 6.684855636158299
$ws.Cells.Item(22, 5).Value = 10.44139176049689
$ws.Cells.Item(22, 7).Value = 3.718576492742545
$ws.Cells.Item(22, 11).Value = 17.64567827062366
$ws.Cells.Item(22, 13).Value = 18.21474902460469
$ws.Cells.Item(22, 14).Value = 26.34480655454643

$ws.Cells.Item(23, 2).Value = 18.49687118531545
$ws.Cells.Item(23, 3).Value = 10.49667584204618
$ws.Cells.Item(23, 4).Value = 6.631578430730693
$ws.Cells.Item(23, 5).Value = 10.34619899109723
$ws.Cells.Item(23, 7).Value = 3.720170539821609
$ws.Cells.Item(23, 11).Value = 17.54620547048389
$ws.Cells.Item(23, 13).Value = 18.12600798538734
$ws.Cells.Item(23, 14).Value = 26.29033180204189

$ws.Cells.Item(24, 2).Value = 17.99146443011321
$ws.Cells.Item(24, 3).Value = 10.1048748618718
$ws.Cells.Item(24, 4).Value = 6.429225960024016
$ws.Cells.Item(24, 5).Value = 9.983328138333885
$ws.Cells.Item(24, 7).Value = 3.726424098226965
$ws.Cells.Item(24, 11).Value = 17.1729397325568
$ws.Cells.Item(24, 13).Value = 17.7949101678032
$ws.Cells.Item(24, 14).Value = 26.08482986767438

$ws.Cells.Item(25, 2).Value = 17.4503361782647
$ws.Cells.Item(25, 3).Value = 9.676870248606726
$ws.Cells.Item(25, 4).Value = 6.211029847245408
$ws.Cells.Item(25, 5).Value = 9.589268149034417
$ws.Cells.Item(25, 7).Value = 3.733636493795135
$ws.Cells.Item(25, 11).Value = 16.77979646640498
$ws.Cells.Item(25, 13).Value = 17.45011394038411
$ws.Cells.Item(25, 14).Value = 25.86523149045736
